$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the BigQuery-style fully-qualified table references: the old
# placeholder project/dataset macros are replaced with the real
# project id / dataset name used in prj-isp-a1777-appl-svil-001.
$ws.Range("A2").Value = "prj-isp-a1777-appl-svil-001.A1777W.VA_A1_HD_ARCHIV_ANAGRA_TITOLI_ON.COD_ABI"
$ws.Range("B2").Value = "prj-isp-a1777-appl-svil-001.A1777W.VA_A1_HD_ARCHIV_ANAGRA_TITOLI_ON.COD_ABI"
$ws.Range("A3").Value = "prj-isp-a1777-appl-svil-001.A1777W.VA_A1_HD_ARCHIV_ANAGRA_TITOLI_ON.NUM_ISTITUTO"
$ws.Range("B3").Value = "prj-isp-a1777-appl-svil-001.A1777W.VA_A1_HD_ARCHIV_ANAGRA_TITOLI_ON.NUM_ISTITUTO"
$ws.Range("A4").Value = "prj-isp-a1777-appl-svil-001.A1777W.VA_A1_HD_ARCHIV_ANAGRA_TITOLI_ON.COD_TITOLO"
$ws.Range("B4").Value = "prj-isp-a1777-appl-svil-001.A1777W.VA_A1_HD_ARCHIV_ANAGRA_TITOLI_ON.COD_TITOLO"
$ws.Range("A5").Value = "prj-isp-a1777-appl-svil-001.A1777W.VA_A1_HD_ARCHIV_ANAGRA_TITOLI_ON.COD_ENTITA"

# Split the merged A:B column width into two independently-sized columns.
$ws.Columns.Item(1).ColumnWidth = 79
$ws.Columns.Item(2).ColumnWidth = 79.1667

# Leave the selection where the author last clicked while working.
$ws.Range("B9").Select()
